# Generate Report for Handback
#
# For both the zh-cn and de-de handback sheets, row 7 (the
# 907d0ec7-bcca-4d2b-a300-42da4df8e172 entry) now has a "Latest Target File"
# hyperlink (column I), a "Latest Target File" name (column J), an updated
# "Latest Handback DateTime" (column K) and an "Error Detail" message
# (column P) describing that the handed-back file was stale. The
# "Error Detail" column (P) is also widened to fit the new text.

$wb = $excel.ActiveWorkbook

function Update-HandbackSheet {
    param(
        [string]$SheetName,
        [string]$TargetFileName,
        [string]$HandbackDateTime,
        [string]$ErrorDetail
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the "Error Detail" column (P / col 16) to fit the new text.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Capture the existing hyperlinks (address/url/display) so they can be
    # rebuilt in order, with the new I7 hyperlink inserted right before A8 -
    # this keeps the relationship ids in the same order as a natural
    # Excel edit (new hyperlink gets rId9, A8's hyperlink shifts to rId10).
    $existing = @()
    foreach ($hl in $ws.Hyperlinks) {
        $existing += , @($hl.Range.Address(), $hl.Address, $hl.TextToDisplay)
    }

    # Find A7's target url - the new I7 hyperlink points at the same
    # ("latest") commit of the markdown file.
    $a7Url = $null
    foreach ($l in $existing) {
        if ($l[0] -eq "`$A`$7") {
            $a7Url = $l[1]
        }
    }

    $ws.Hyperlinks.Delete()

    foreach ($l in $existing) {
        $addr = $l[0]
        $url = $l[1]
        $disp = $l[2]
        if ($addr -eq "`$A`$8") {
            $ws.Hyperlinks.Add($ws.Range("I7"), $a7Url, "", "", $TargetFileName)
            $ws.Range("I7").Font.Underline = 2
            $ws.Range("I7").Font.Color = 15570276
        }
        $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $disp)
    }

    $ws.Range("J7").Value = $TargetFileName
    $ws.Range("K7").Value = $HandbackDateTime
    $ws.Range("P7").Value = $ErrorDetail
}

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c979792057a57f91648e92b47e45096549ffc4e2/e2e/907d0ec7-bcca-4d2b-a300-42da4df8e172.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd039702c0ad32bbb44c87f955b42a6ad141a1d2/e2e/907d0ec7-bcca-4d2b-a300-42da4df8e172.md."

Update-HandbackSheet "zh-cn" `
    "907d0ec7-bcca-4d2b-a300-42da4df8e172.c3fd4bb07494907343e7002cbbba9e731b157864.zh-cn.xlf" `
    "2016-08-31 02:49:34" `
    $errorDetail

Update-HandbackSheet "de-de" `
    "907d0ec7-bcca-4d2b-a300-42da4df8e172.c3fd4bb07494907343e7002cbbba9e731b157864.de-de.xlf" `
    "2016-08-31 02:49:40" `
    $errorDetail
